$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row renames
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Standardize "de" -> "De" capitalization in various place names
$ws.Range("B5").Value = "San Cristóbal De Las Casas"
$ws.Range("A13").Value = "Estado De México"
$ws.Range("A15").Value = "Guanajuato"
$ws.Range("B15").Value = "Apaseo El Alto"
$ws.Range("B17").Value = "Acapulco De Juárez"
$ws.Range("B20").Value = "Mártir De Cuilapan"
$ws.Range("B33").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B34").Value = "Santo Domingo De Morelos"
$ws.Range("B40").Value = "Amealco De Bonfil"
$ws.Range("B42").Value = "Landa De Matamoros"
$ws.Range("B51").Value = "Hueyapan De Ocampo"
$ws.Range("B55").Value = "Noria De Ángeles"

# Minor float precision fix
$ws.Range("D56").Value = 0.09677419354838708

# Delete footer rows (60-64): sample size, source, elaborated by, secretariat, date
$ws.Range("A60:D64").EntireRow.Delete()
